$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new row (row 8) with the next day's gold data, matching the
# formatting pattern already used by the preceding rows (6 and 7).
$ws.Range("A8").Value = "29-09-2025"
$ws.Range("B8").Value = "The price of gold in India today is ₹11,640 per gram for 24 karat gold, ₹10,670 per gram for 22 karat gold and ₹8,730 per gram for 18 karat gold (also called 999 gold)."

# Copy the formatting from the row above so the new row matches style-wise.
$ws.Range("A7:B7").Copy()
$ws.Range("A8:B8").PasteSpecial(-4122)
